# Fruta / hortaliza, semanal
# Insert 3 new weekly records at the top of the Espárragos data block
# (rows 73-75), pushing the existing rows 73-96 down to 76-99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 73; this shifts the
# existing rows 73:96 down to 76:99, preserving all of their data/styles.
$ws.Rows("73:75").Insert()

# Row 73 - new record (Extra, Región del Maule)
$ws.Range("A73").Value = 10
$ws.Range("B73").Value = "Vega Modelo de Temuco"
$ws.Range("C73").Value = "La Araucanía"
$ws.Range("D73").Value = 45215
$ws.Range("E73").Value = 9
$ws.Range("F73").Value = 300000000
$ws.Range("G73").Value = "Espárragos"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Extra"
$ws.Range("J73").Value = 200
$ws.Range("K73").Value = 2000
$ws.Range("L73").Value = 2000
$ws.Range("M73").Value = 2000
$ws.Range("N73").Value = "$/kilo"
$ws.Range("O73").Value = "Región del Maule"
$ws.Range("P73").Value = 2000
$ws.Range("Q73").Value = 1
$ws.Range("R73").Value = "Hortaliza"

# Row 74 - new record (Primera, Región de La Araucanía)
$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 45215
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = 300000000
$ws.Range("G74").Value = "Espárragos"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 800
$ws.Range("K74").Value = 1600
$ws.Range("L74").Value = 1600
$ws.Range("M74").Value = 1600
$ws.Range("N74").Value = "$/kilo"
$ws.Range("O74").Value = "Región de La Araucanía"
$ws.Range("P74").Value = 1600
$ws.Range("Q74").Value = 1
$ws.Range("R74").Value = "Hortaliza"

# Row 75 - new record (Primera, Región del Maule)
$ws.Range("A75").Value = 10
$ws.Range("B75").Value = "Vega Modelo de Temuco"
$ws.Range("C75").Value = "La Araucanía"
$ws.Range("D75").Value = 45215
$ws.Range("E75").Value = 9
$ws.Range("F75").Value = 300000000
$ws.Range("G75").Value = "Espárragos"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 300
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 1500
$ws.Range("M75").Value = 1500
$ws.Range("N75").Value = "$/kilo"
$ws.Range("O75").Value = "Región del Maule"
$ws.Range("P75").Value = 1500
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = "Hortaliza"
